# Update "paises.xlsx" (Pais sheet) with newer COVID country stats
# and bump the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" banner in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 08:03"

# Row 27 - Ucrania
$ws.Range("B27").Value = 143030
$ws.Range("C27").Value = 2551
$ws.Range("E27").Value = 76505
$ws.Range("G27").Value = 45
$ws.Range("H27").Value = 2979

# Row 28 - Israel
$ws.Range("B28").Value = 138719
$ws.Range("C28").Value = 1154
$ws.Range("D28").Value = 107600
$ws.Range("E28").Value = 30079

# Row 33 - Kazajistan
$ws.Range("D33").Value = 100042
$ws.Range("E33").Value = 4822

# Row 63 - Kirguistan
$ws.Range("B63").Value = 44613
$ws.Range("C63").Value = 87
$ws.Range("D63").Value = 40336
$ws.Range("E63").Value = 3216
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 1061

# Row 64 - Uzbekistan
$ws.Range("B64").Value = 44557
$ws.Range("C64").Value = 276
$ws.Range("D64").Value = 41898
$ws.Range("E64").Value = 2297
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 362

# Row 202 - Fiyi
$ws.Range("B202").Value = 32
$ws.Range("C202").Value = 1
$ws.Range("E202").Value = 6
